$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: push the existing 2021-Q3 summary row down to row 3 and
#    write a brand-new 2022-Q3 summary row into row 2 (the new quarter is
#    reported first).
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Clone row 2's formatting onto the (new) row 3 before touching any values,
# so A3 picks up the same style (s="2") that A2 already has.
$totals.Range("A2").Copy()
$totals.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2021-Q3"
$totals.Range("C3").Value = 2
$totals.Range("D3").Value = 0.12

$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0.09

# ---------------------------------------------------------------------------
# 2) Insert a brand-new worksheet named "2022-Q3" immediately before the
#    existing "2021-Q3" sheet. Because "2021-Q3" is currently sheet2.xml,
#    this shifts it to sheet3.xml while the new sheet becomes sheet2.xml -
#    matching the upstream diff's file layout (and sheet order).
# ---------------------------------------------------------------------------
$q3_2021 = $wb.Worksheets.Item("2021-Q3")
$q3_2022 = $wb.Worksheets.Add($q3_2021)
$q3_2022.Name = "2022-Q3"

# Match the page margins used by every other sheet in the workbook
# (0.75in/0.75in/1in/1in/0.5in/0.5in).
$q3_2022.PageSetup.LeftMargin = 54
$q3_2022.PageSetup.RightMargin = 54
$q3_2022.PageSetup.TopMargin = 72
$q3_2022.PageSetup.BottomMargin = 72
$q3_2022.PageSetup.HeaderMargin = 36
$q3_2022.PageSetup.FooterMargin = 36

# Header row (B1:H1) - reuse the same bold/centered/bordered style (s="2")
# that "总计" uses for its own header row.
$q3_2022.Range("B1").Value = "基金代码"
$q3_2022.Range("C1").Value = "基金名称"
$q3_2022.Range("D1").Value = "基金规模"
$q3_2022.Range("E1").Value = "股票总仓位"
$q3_2022.Range("F1").Value = "仓位占比"
$q3_2022.Range("G1").Value = "持有市值(亿元)"
$q3_2022.Range("H1").Value = "仓位排名"

$totals.Range("B1").Copy()
$q3_2022.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

$totals.Range("A2").Copy()
$q3_2022.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# Data row 2 - the fund code / numeric-looking figures must stay TEXT
# (leading zero, fixed decimals matter) so force a Text number format while
# writing them, then drop back to the default "Normal" style so no stray
# explicit style ends up on the cell (matches the unstyled <c> cells in the
# upstream sheet).
$q3_2022.Range("A2").Value = 0

$textCells = $q3_2022.Range("B2:G2")
$textCells.NumberFormat = "@"
$q3_2022.Range("B2").Value = "010204"
$q3_2022.Range("C2").Value = "中银港股通优势成长股票"
$q3_2022.Range("D2").Value = "2.76"
$q3_2022.Range("E2").Value = "79.66"
$q3_2022.Range("F2").Value = "3.40"
$q3_2022.Range("G2").Value = "0.0938"
$textCells.Style = "Normal"

$q3_2022.Range("H2").Value = 8

$totals.Activate()
$totals.Range("A1").Select()
